$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits ---

# Defaults row: Spring constant (k), C3: 10 -> 10000
$ws.Range("C3").Value = 10000

# Values row (row 4):
# mega_arrays (G4): False -> True  (copy existing "True" text from G3 so the
# shared-string gets reused instead of Excel auto-typing it as a boolean)
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4163)

# GPU Compute (L4): True -> False (copy existing "False" text from P3)
$ws.Range("P3").Copy()
$ws.Range("L4").PasteSpecial(-4163)

# Custom Shape? (P4): True -> False (copy existing "False" text from P3)
$ws.Range("P3").Copy()
$ws.Range("P4").PasteSpecial(-4163)

# Repeats (I4): 1 -> 150
$ws.Range("I4").Value = "150"

$excel.CutCopyMode = $false

# --- View / selection edits ---
$ws.Range("I4").Select()

$wb.Save()
